# Correcciones en reglas del documento stock actual
# Updates "Diferencia Stock" (L), recomputed sales/profit targets (M/N),
# real stock (P), corrected order (Q) and final order (U) for several
# articles, re-sorts/relabels the BONSAI size rows (35-38), hides rows
# whose corrected final order dropped to 0, and refreshes the summary
# metrics block at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple "Diferencia Stock" (column L) corrections ---------------------
$ws.Range("L3").Value = 1
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("L11").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("L14").Value = -1
$ws.Range("L18").Value = -1
$ws.Range("L19").Value = -1
$ws.Range("L24").Value = 1
$ws.Range("L32").Value = -1
$ws.Range("L34").Value = 1
$ws.Range("L40").Value = -1
$ws.Range("L41").Value = 1
$ws.Range("L46").Value = 1
$ws.Range("L58").Value = 1
$ws.Range("L60").Value = 1
$ws.Range("L61").Value = 1
$ws.Range("L63").Value = 1

# --- Row 31 (DIEFFENBACHIA CAMILA) now covered by real stock --------------
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("P31").Value = 7
$ws.Range("Q31").Value = 0
$ws.Range("U31").Value = 0
$ws.Rows.Item(31).Hidden = $true

# --- Row 33 (ANTHURIUM ANDRAEANUM MIX) now covered by real stock ----------
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("P33").Value = 15
$ws.Range("Q33").Value = 0
$ws.Range("U33").Value = 0
$ws.Rows.Item(33).Hidden = $true

# --- BONSAI sizes (rows 35-38) re-sorted by size and refreshed -------------
# New order: 5A, 6A, 8A, 7ASUP (was 7ASUP, 8A, 5A, 6A)
$ws.Range("C35").Value = "5A       "
$ws.Range("G35").Value = 18.12
$ws.Range("H35").Value = 7.25
$ws.Range("J35").Value = "REDUCIR 15%"
$ws.Range("L35").Value = -1
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("P35").Value = 12
$ws.Range("Q35").Value = 0
$ws.Range("S35").Value = 3
$ws.Range("T35").Value = 3
$ws.Range("U35").Value = 0

$ws.Range("C36").Value = "6A       "
$ws.Range("G36").Value = 29.62
$ws.Range("H36").Value = 11.85
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("P36").Value = 10
$ws.Range("Q36").Value = 0
$ws.Range("U36").Value = 0

$ws.Range("C37").Value = "8A       "
$ws.Range("G37").Value = 48.92
$ws.Range("H37").Value = 19.57
$ws.Range("J37").Value = "REDUCIR 50%"
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("P37").Value = 10
$ws.Range("Q37").Value = 0
$ws.Range("S37").Value = 0
$ws.Range("T37").Value = 0
$ws.Range("U37").Value = 0

$ws.Range("C38").Value = "7ASUP    "
$ws.Range("G38").Value = 49.42
$ws.Range("H38").Value = 19.77
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("P38").Value = 2
$ws.Range("Q38").Value = 0
$ws.Range("U38").Value = 0

$ws.Rows.Item(35).Hidden = $true
$ws.Rows.Item(36).Hidden = $true
$ws.Rows.Item(37).Hidden = $true
$ws.Rows.Item(38).Hidden = $true

# --- Row 47 (NEPHROLEPIS BOSTON) now covered by real stock -----------------
$ws.Range("L47").Value = 1
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("P47").Value = 18
$ws.Range("Q47").Value = 0
$ws.Range("U47").Value = 0
$ws.Rows.Item(47).Hidden = $true

# --- Row 48 (MONSTERA) now covered by real stock ----------------------------
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("P48").Value = 8
$ws.Range("Q48").Value = 0
$ws.Range("U48").Value = 0
$ws.Rows.Item(48).Hidden = $true

# --- Row 55 (SENECIO MIX) now covered by real stock -------------------------
$ws.Range("M55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("P55").Value = 4
$ws.Range("Q55").Value = 0
$ws.Range("U55").Value = 0
$ws.Rows.Item(55).Hidden = $true

# --- Row 56 (PLANTA CRASA M85A15) partially covered by real stock ----------
$ws.Range("L56").Value = -3
$ws.Range("M56").Value = 30.8
$ws.Range("N56").Value = 18.48
$ws.Range("P56").Value = 20
$ws.Range("Q56").Value = 13
$ws.Range("U56").Value = 14

# --- Row 57 (PLANTA CRASA M55A10) now covered by real stock -----------------
$ws.Range("L57").Value = 1
$ws.Range("M57").Value = 0
$ws.Range("N57").Value = 0
$ws.Range("P57").Value = 108
$ws.Range("Q57").Value = 0
$ws.Range("U57").Value = 0
$ws.Rows.Item(57).Hidden = $true

# --- Refresh the summary metrics block --------------------------------------
$ws.Range("C66").Value = 153
$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value = "1690.32€"
$ws.Range("C77").Value = 7
